# Game-idea document update
# Applies the "docs: update Game Idea" edits:
#  - Rename the game / change the tagline
#  - Change genre text
#  - Update the controls list (keys + actions)
#  - Shrink the "Tổng quan" (overview) paragraph down to just its label
#  - Remove the "Ví dụ minh họa" bullet and the screenshot that followed it
#  - Move the _GoBack bookmark onto the end of the controls section

$d = $word.ActiveDocument

# --- 1. Title line (tagline) ---------------------------------------------
# "Tên : INTERN 'video' game"  ->  "Tên : Catgirl & Shotgun Project"
$d.Content.Find.Execute(
    " : INTERN 'video' game", $true, $false, $false, $false, $false,
    $true, 1, $false, " : Catgirl & Shotgun Project", 2) | Out-Null

# --- 2. Genre line ----------------------------------------------------------
# "Thể loại" stays bold; only the non-bold remainder changes.
# ": Puzzle Platformer" -> ": 2D Top-down Shooter"
$d.Content.Find.Execute(
    ": Puzzle Platformer", $true, $false, $false, $false, $false,
    $true, 1, $false, ": 2D Top-down Shooter", 2) | Out-Null

# --- 3. Controls list -------------------------------------------------------
# Movement keys gain the Arrow-Keys alternative.
$d.Content.Find.Execute(
    "[WASD]: Di chuyển nhân vật.", $true, $false, $false, $false, $false,
    $true, 1, $false, "[WASD]/[Arrow Keys]: Di chuyển nhân vật.", 2) | Out-Null

# [Space]: Nhảy.  ->  [Q]: Sử dụng kĩ năng
$d.Content.Find.Execute(
    "[Space]: Nhảy.", $true, $false, $false, $false, $false,
    $true, 1, $false, "[Q]: Sử dụng kĩ năng", 2) | Out-Null

# [E]: Tương tác với các đối tượng trong game.  ->  [Left click]/[Space]: Bắn
$d.Content.Find.Execute(
    "[E]: Tương tác với các đối tượng trong game.", $true, $false, $false, $false, $false,
    $true, 1, $false, "[Left click]/[Space]: Bắn", 2) | Out-Null

# --- 4. Move the _GoBack bookmark to right after "Bắn" ----------------------
# (It used to sit on the "Ví dụ minh họa" bullet, which is being removed.)
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") { $bm.Delete() }
}

$controlsPara = $d.Paragraphs(8)
# Add a throwaway marker character so we can Find it and Collapse cleanly
# onto the real end-of-text position (collapsing straight onto the
# paragraph-mark boundary lands the bookmark in the wrong paragraph).
$controlsPara.Range.InsertAfter("#")
$marker = $controlsPara.Range.Duplicate
$marker.Find.Execute("#", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$marker.Collapse(1)
$d.Bookmarks.Add("_GoBack", $marker) | Out-Null
$markerChar = $d.Range($marker.Start, $marker.Start + 1)
$markerChar.Delete()

# --- 5. Shorten the "Tổng quan" paragraph -----------------------------------
$d.Content.Find.Execute(
    "Tổng quan: *đố.", $false, $false, $true, $false, $false,
    $true, 1, $false, "Tổng quan: ", 2) | Out-Null

# --- 6. Remove the "Ví dụ minh họa" bullet and the screenshot paragraph ----
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt.StartsWith("Ví dụ minh họa")) {
        $exampleIndex = $i
        break
    }
}

if ($exampleIndex) {
    $examplePara = $d.Paragraphs($exampleIndex)
    $examplePara.Range.Delete()

    # The picture now lives in what has become paragraph $exampleIndex.
    while ($d.Shapes.Count -gt 0) {
        $d.Shapes(1).Delete()
    }
    $pictureParaIndex = $exampleIndex
    $picturePara = $d.Paragraphs($pictureParaIndex)
    $picturePara.Range.Delete()
}
